# The deck currently ships two theme parts that are, in effect, swapped:
#   - theme2.xml (the theme actually applied to the Slide Master / every
#     slide) holds the "Integral" green/olive colour scheme.
#   - theme1.xml (only wired up via the Notes Master relationship) still
#     holds the stock blue "Office Theme" colour scheme.
# The authored change flips which physical part carries which palette, so
# the Slide Master's theme becomes the plain "Office Theme" colours (and
# the Notes Master side picks up "Integral"). The only things that
# actually differ between the two theme parts are the 12 theme colours
# (plus their display names) - font scheme and format scheme are already
# byte-identical - so we reproduce the visible effect by rewriting the
# Slide Master's theme colour scheme, one swatch at a time, through the
# standard Theme / ThemeColorScheme object model.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$tcs = $theme.ThemeColorScheme

# Office Theme palette (RGB() takes R + G*256 + B*65536, i.e. 0xBBGGRR).
$tcs.Item(1).RGB  = 0x000000   # dk1      -> 000000
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 0x6A5444   # dk2      -> 44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # accent2  -> ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$tcs.Item(9).RGB  = 0xC47244   # accent5  -> 4472C4
$tcs.Item(10).RGB = 0x47AD70   # accent6  -> 70AD47
$tcs.Item(11).RGB = 0xC16305   # hlink    -> 0563C1
$tcs.Item(12).RGB = 0x724F95   # folHlink -> 954F72
